$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume% (E) columns store numeric-looking values as text
# (e.g. "159.00", "1.000", "30.135.52"). Force text format before writing
# so Excel does not auto-convert them to numbers and lose formatting,
# then clear the explicit format again so no new cell style is introduced.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.135.52"
$ws.Range("E2").Value = "  -4.25%  "

$ws.Range("D3").Value = "1.910.39"
$ws.Range("E3").Value = "  -3.91%  "

$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "245.18"
$ws.Range("E5").Value = "  -3.18%  "

$ws.Range("D6").Value = "0.7085"
$ws.Range("E6").Value = "  -13.53%  "

$ws.Range("D7").Value = "0.9975"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("D8").Value = "0.3224"
$ws.Range("E8").Value = "  -5.30%  "

$ws.Range("D9").Value = "26.05"
$ws.Range("E9").Value = "  +1.54%  "

$ws.Range("D10").Value = "0.06822"
$ws.Range("E10").Value = "  -2.64%  "

$ws.Range("D11").Value = "0.7854"
$ws.Range("E11").Value = "  -6.81%  "

$ws.Range("D12").Value = "0.07911"
$ws.Range("E12").Value = "  -2.45%  "

$ws.Range("D13").Value = "1.909.60"
$ws.Range("E13").Value = "  -3.97%  "

$ws.Range("D14").Value = "5.364"
$ws.Range("E14").Value = "  -2.12%  "

$ws.Range("D15").Value = "93.62"
$ws.Range("E15").Value = "  -7.99%  "

$ws.Range("D16").Value = "14.34"
$ws.Range("E16").Value = "  +2.47%  "

$ws.Range("D17").Value = "258.10"
$ws.Range("E17").Value = "  -5.49%  "

$ws.Range("D18").Value = "30.147.86"
$ws.Range("E18").Value = "  -4.17%  "

$ws.Range("D19").Value = "5.776"
$ws.Range("E19").Value = "  +0.85%  "

$ws.Range("D20").Value = "0.000007820"
$ws.Range("E20").Value = "  -2.42%  "

$ws.Range("D21").Value = "2.163.53"
$ws.Range("E21").Value = "  -3.52%  "

$ws.Range("D22").Value = "0.9984"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").Value = "6.787"
$ws.Range("E24").Value = "  -2.30%  "

$ws.Range("D25").Value = "9.550"
$ws.Range("E25").Value = "  -1.34%  "

$ws.Range("D26").Value = "159.00"
$ws.Range("E26").Value = "  -3.90%  "

$ws.Range("D27").Value = "0.1311"
$ws.Range("E27").Value = "  -16.53%  "

$ws.Range("D28").Value = "18.66"
$ws.Range("E28").Value = "  -5.45%  "

$ws.Range("D29").Value = "2.191"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").Value = "1.359"
$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").Value = "1.538"
$ws.Range("E31").Value = "  -1.85%  "

$ws.Range("D32").Value = "4.387"
$ws.Range("E32").Value = "  -3.84%  "

$ws.Range("D33").Value = "4.159"
$ws.Range("E33").Value = "  -4.14%  "

$ws.Range("D34").Value = "0.04993"
$ws.Range("E34").Value = "  -4.02%  "

$ws.Range("D35").Value = "1.179"
$ws.Range("E35").Value = "  -2.91%  "

$ws.Range("D36").Value = "0.7371"
$ws.Range("E36").Value = "  -1.82%  "

$ws.Range("D37").Value = "2.715"
$ws.Range("E37").Value = "  -3.13%  "

$ws.Range("D38").Value = "0.01921"
$ws.Range("E38").Value = "  -4.05%  "

$ws.Range("D39").Value = "2.789"
$ws.Range("E39").Value = "  -4.62%  "

$ws.Range("D40").Value = "78.98"
$ws.Range("E40").Value = "  +0.33%  "

$ws.Range("D41").Value = "6.408"
$ws.Range("E41").Value = "  -3.40%  "

$ws.Range("D42").Value = "0.4383"
$ws.Range("E42").Value = "  -6.12%  "

$ws.Range("D43").Value = "2.001"
$ws.Range("E43").Value = "  -3.08%  "

$ws.Range("D44").Value = "0.9972"
$ws.Range("E44").Value = "  -0.33%  "

$ws.Range("D45").Value = "0.8298"
$ws.Range("E45").Value = "  -2.83%  "

$ws.Range("D46").Value = "101.64"
$ws.Range("E46").Value = "  -4.40%  "

$ws.Range("D47").Value = "9.591"
$ws.Range("E47").Value = "  -4.27%  "

$ws.Range("D48").Value = "7.194"
$ws.Range("E48").Value = "  -3.97%  "

$ws.Range("D49").Value = "35.77"
$ws.Range("E49").Value = "  -2.17%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05898"
$ws.Range("E50").Value = "  -1.33%  "

$ws.Range("D51").Value = "1.462"
$ws.Range("E51").Value = "  +1.55%  "

# Remove the explicit text-format style again so the cells keep their
# original (default/general) style, matching the source formatting.
$priceRange.ClearFormats()
